$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.875.01"
$ws.Range("E2").Value = "  +0.76%  "

# Row 3
$ws.Range("D3").Value = "2.533.65"
$ws.Range("E3").Value = "  +0.71%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'315.96"
$ws.Range("E5").Value = "  +0.61%  "

# Row 6
$ws.Range("D6").Value = "'96.09"
$ws.Range("E6").Value = "  +0.39%  "

# Row 7
$ws.Range("E7").Value = "  -1.45%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  -0.91%  "

# Row 10
$ws.Range("D10").Value = "'36.19"
$ws.Range("E10").Value = "  -0.48%  "

# Row 12
$ws.Range("D12").Value = "'7.57"
$ws.Range("E12").Value = "  -0.56%  "

# Row 13
$ws.Range("E13").Value = "  -3.47%  "

# Row 14
$ws.Range("D14").Value = "2.922.14"
$ws.Range("E14").Value = "  +0.73%  "

# Row 15
$ws.Range("D15").Value = "2.554.50"
$ws.Range("E15").Value = "  +1.91%  "

# Row 16
$ws.Range("D16").Value = "'15.25"
$ws.Range("E16").Value = "  -2.21%  "

# Row 17
$ws.Range("D17").Value = "'0.856"
$ws.Range("E17").Value = "  -0.16%  "

# Row 18
$ws.Range("D18").Value = "42.930.47"
$ws.Range("E18").Value = "  +0.92%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'6.79"
$ws.Range("E19").Value = "  +4.44%  "

# Row 20
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.88"
$ws.Range("E20").Value = "  -0.23%  "

# Row 21
$ws.Range("E21").Value = "  -0.98%  "

# Row 22
$ws.Range("D22").Value = "'69.96"
$ws.Range("E22").Value = "  -2.19%  "

# Row 23
$ws.Range("D23").Value = "'254.14"
$ws.Range("E23").Value = "  +0.29%  "

# Row 24
$ws.Range("D24").Value = "'2.95"
$ws.Range("E24").Value = "  -1.27%  "

# Row 25
$ws.Range("E25").Value = "  +2.28%  "

# Row 26
$ws.Range("D26").Value = "'26.79"
$ws.Range("E26").Value = "  -0.76%  "

# Row 27
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("D28").Value = "'2.42"
$ws.Range("E28").Value = "  +3.29%  "

# Row 29
$ws.Range("D29").Value = "'40.78"
$ws.Range("E29").Value = "  +8.03%  "

# Row 30
$ws.Range("D30").Value = "'10.38"
$ws.Range("E30").Value = "  +2.49%  "

# Row 31
$ws.Range("D31").Value = "'5.92"
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$ws.Range("E32").Value = "  +2.04%  "

# Row 33
$ws.Range("D33").Value = "'2.19"
$ws.Range("E33").Value = "  +5.41%  "

# Row 34
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'3.35"
$ws.Range("E34").Value = "  +1.22%  "

# Row 35
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'19.22"
$ws.Range("E35").Value = "  +0.32%  "

# Row 36
$ws.Range("E36").Value = "  +2.06%  "

# Row 37
$ws.Range("D37").Value = "'0.0782"
$ws.Range("E37").Value = "  -0.50%  "

# Row 38
$ws.Range("E38").Value = "  -1.03%  "

# Row 39
$ws.Range("E39").Value = "  -0.92%  "

# Row 40
$ws.Range("D40").Value = "'23.49"
$ws.Range("E40").Value = "  -3.63%  "

# Row 41
$ws.Range("E41").Value = "  +15.03%  "

# Row 42
$ws.Range("D42").Value = "'3.86"
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("D43").Value = "'0.0304"
$ws.Range("E43").Value = "  +0.58%  "

# Row 45
$ws.Range("E45").Value = "  -1.62%  "

# Row 46
$ws.Range("D46").Value = "2.042.30"
$ws.Range("E46").Value = "  +0.98%  "

# Row 47
$ws.Range("D47").Value = "'85.03"
$ws.Range("E47").Value = "  +0.71%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'9.03"
$ws.Range("E48").Value = "  +1.05%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'107.08"
$ws.Range("E49").Value = "  +5.78%  "

# Row 50
$ws.Range("D50").Value = "'75.26"
$ws.Range("E50").Value = "  +2.95%  "

# Row 51
$ws.Range("D51").Value = "2.777.18"
$ws.Range("E51").Value = "  +0.73%  "
